# Agregue mi estimacion de tareas. Cristian
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimación de tareas")

$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 10
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 8

$ws.Activate()
$ws.Range("G15").Select()
